$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 165852
$ws.Range("C4").Value = 156783
$ws.Range("C5").Value = 9069
$ws.Range("C8").Value = 65.23999999999999
